{"js": "// Replace each three-digit-division answer cell's text with its new value.\n// Each \"before\" text is unique in the document, so a direct search/replace\n// per pair is safe and unambiguous.\nconst replacements = [\n  [\"578\u00f78=72, 2\", \"959\u00f73=319, 2\"],\n  [\"501\u00f73=167, 0\", \"404\u00f78=50, 4\"],\n  [\"276\u00f77=39, 3\", \"413\u00f77=59, 0\"],\n  [\"454\u00f77=64, 6\", \"189\u00f76=31, 3\"],\n  [\"221\u00f75=44, 1\", \"700\u00f74=175, 0\"],\n  [\"778\u00f74=194, 2\", \"370\u00f77=52, 6\"],\n  [\"100\u00f73=33, 1\", \"942\u00f77=134, 4\"],\n  [\"235\u00f77=33, 4\", \"506\u00f73=168, 2\"],\n  [\"397\u00f78=49, 5\", \"221\u00f78=27, 5\"],\n  [\"123\u00f79=13, 6\", \"974\u00f75=194, 4\"],\n  [\"892\u00f77=127, 3\", \"357\u00f76=59, 3\"],\n  [\"893\u00f78=111, 5\", \"965\u00f75=193, 0\"],\n  [\"120\u00f74=30, 0\", \"642\u00f75=128, 2\"],\n  [\"185\u00f73=61, 2\", \"994\u00f76=165, 4\"],\n  [\"685\u00f78=85, 5\", \"174\u00f73=58, 0\"],\n  [\"491\u00f72=245, 1\", \"452\u00f73=150, 2\"],\n  [\"589\u00f76=98, 1\", \"631\u00f72=315, 1\"],\n  [\"240\u00f74=60, 0\", \"670\u00f77=95, 5\"],\n  [\"773\u00f73=257, 2\", \"721\u00f78=90, 1\"],\n  [\"602\u00f75=120, 2\", \"881\u00f74=220, 1\"],\n  [\"733\u00f73=244, 1\", \"533\u00f77=76, 1\"],\n  [\"870\u00f78=108, 6\", \"881\u00f76=146, 5\"],\n  [\"428\u00f76=71, 2\", \"877\u00f73=292, 1\"],\n  [\"202\u00f78=25, 2\", \"209\u00f78=26, 1\"],\n  [\"745\u00f72=372, 1\", \"185\u00f72=92, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit-division answer cell's text with its new value.\n# Each \"before\" text occurs exactly once in the document, so Find/Replace\n# (wdReplaceAll) per pair is unambiguous and safe.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"578\u00f78=72, 2\", \"959\u00f73=319, 2\"),\n    @(\"501\u00f73=167, 0\", \"404\u00f78=50, 4\"),\n    @(\"276\u00f77=39, 3\", \"413\u00f77=59, 0\"),\n    @(\"454\u00f77=64, 6\", \"189\u00f76=31, 3\"),\n    @(\"221\u00f75=44, 1\", \"700\u00f74=175, 0\"),\n    @(\"778\u00f74=194, 2\", \"370\u00f77=52, 6\"),\n    @(\"100\u00f73=33, 1\", \"942\u00f77=134, 4\"),\n    @(\"235\u00f77=33, 4\", \"506\u00f73=168, 2\"),\n    @(\"397\u00f78=49, 5\", \"221\u00f78=27, 5\"),\n    @(\"123\u00f79=13, 6\", \"974\u00f75=194, 4\"),\n    @(\"892\u00f77=127, 3\", \"357\u00f76=59, 3\"),\n    @(\"893\u00f78=111, 5\", \"965\u00f75=193, 0\"),\n    @(\"120\u00f74=30, 0\", \"642\u00f75=128, 2\"),\n    @(\"185\u00f73=61, 2\", \"994\u00f76=165, 4\"),\n    @(\"685\u00f78=85, 5\", \"174\u00f73=58, 0\"),\n    @(\"491\u00f72=245, 1\", \"452\u00f73=150, 2\"),\n    @(\"589\u00f76=98, 1\", \"631\u00f72=315, 1\"),\n    @(\"240\u00f74=60, 0\", \"670\u00f77=95, 5\"),\n    @(\"773\u00f73=257, 2\", \"721\u00f78=90, 1\"),\n    @(\"602\u00f75=120, 2\", \"881\u00f74=220, 1\"),\n    @(\"733\u00f73=244, 1\", \"533\u00f77=76, 1\"),\n    @(\"870\u00f78=108, 6\", \"881\u00f76=146, 5\"),\n    @(\"428\u00f76=71, 2\", \"877\u00f73=292, 1\"),\n    @(\"202\u00f78=25, 2\", \"209\u00f78=26, 1\"),\n    @(\"745\u00f72=372, 1\", \"185\u00f72=92, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
